$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 - 2020-03-31 (serial 43921)
$ws.Cells.Item(34, 1).Value = 43921
$ws.Cells.Item(34, 1).NumberFormat = $ws.Cells.Item(33, 1).NumberFormat
$ws.Cells.Item(34, 2).Value = 48
$ws.Cells.Item(34, 3).Value = 600
$ws.Cells.Item(34, 4).Value = 10
$ws.Cells.Item(34, 5).Value = 47
$ws.Cells.Item(34, 6).Value = 58
$ws.Cells.Item(34, 7).Value = 647
$ws.Cells.Item(34, 8).Value = 11
$ws.Cells.Item(34, 9).Value = 74
$ws.Cells.Item(34, 10).Value = 14
$ws.Cells.Item(34, 11).Value = 30
$ws.Cells.Item(34, 12).Value = 2
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(34, 14).Value = 1
$ws.Cells.Item(34, 15).Value = 343
$ws.Cells.Item(34, 16).Value = 188
$ws.Cells.Item(34, 17).Value = 110
$ws.Cells.Item(34, 18).Value = 6

# Row 35 - 2020-04-01 (serial 43922)
$ws.Cells.Item(35, 1).Value = 43922
$ws.Cells.Item(35, 1).NumberFormat = $ws.Cells.Item(33, 1).NumberFormat
$ws.Cells.Item(35, 2).Value = 47
$ws.Cells.Item(35, 3).Value = 647
$ws.Cells.Item(35, 4).Value = 14
$ws.Cells.Item(35, 5).Value = 61
$ws.Cells.Item(35, 6).Value = 61
$ws.Cells.Item(35, 7).Value = 708
$ws.Cells.Item(35, 8).Value = 8
$ws.Cells.Item(35, 9).Value = 82
$ws.Cells.Item(35, 10).Value = 16
$ws.Cells.Item(35, 11).Value = 32
$ws.Cells.Item(35, 12).Value = 2
$ws.Cells.Item(35, 13).Value = 0
$ws.Cells.Item(35, 14).Value = 1
